$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the wrapper element used to embed the MODS datastream update.
# Previously cells C2/W2 wrapped the MODS payload using the legacy
# <update type="MODS">...</update> element; now they use the
# <datastream type="md_descriptive" operation="update">...</datastream> element.
$ws.Range("C2").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">'
$ws.Range("W2").Value = '</mods:mods></datastream></object>'
